# DataSource - Pago.xlsx: fix URL to point to BillingCenter (BC) instead of
# PolicyCenter (PC), so the Ranorex data source works against BC.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newUrl = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/bc/BillingCenter.do"

# Column C holds the URL used by every row (rows 2-53).
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 3).Value = $newUrl
}

# Reflect the cursor/selection left by the edit (column C was selected).
$ws.Range("C2:C53").Select()
